# Scheduled-runner price refresh: update market-board derived columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) for a batch of leve rows
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets. Some profit cells had
# no valid value this run (e.g. division produced nothing), so those are
# cleared rather than written.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 37041468
$ws.Range("I62").Value = 100002740
$ws.Range("J62").Value = 5431.294
$ws.Range("K62").Value = 100002740
$ws.Range("L62").Value = 5431.294
$ws.Range("M62").Value = -100002116
$ws.Range("N62").Value = -6679.294

$ws.Range("H65").Value = 37041468
$ws.Range("I65").Value = 100002740
$ws.Range("J65").Value = 5431.294
$ws.Range("K65").Value = 500013700
$ws.Range("L65").Value = 27156.47
$ws.Range("M65").Value = -500010580
$ws.Range("N65").Value = -33396.47

$ws.Range("H98").Value = 1633.1875
$ws.Range("I98").Value = 1587.0834
$ws.Range("J98").Value = 1771.5
$ws.Range("K98").Value = 1587.0834
$ws.Range("L98").Value = 1771.5
$ws.Range("M98").Value = -89.08339999999998
$ws.Range("N98").Value = -4767.5

$ws.Range("H100").Value = 3580.5
$ws.Range("I100").Value = 3400.7144
$ws.Range("J100").Value = 4000
$ws.Range("K100").Value = 3400.7144
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -2859.7144
$ws.Range("N100").Value = -5082

$ws.Range("H122").Value = 1633.1875
$ws.Range("I122").Value = 1587.0834
$ws.Range("J122").Value = 1771.5
$ws.Range("K122").Value = 4761.2502
$ws.Range("L122").Value = 5314.5
$ws.Range("M122").Value = -2311.2502
$ws.Range("N122").Value = -10214.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 90909090
$ws.Range("I45").Value = 90909090
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 90909090
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -90908713
$ws.Range("N45").ClearContents()

$ws.Range("H102").Value = 2517
$ws.Range("I102").Value = 2293.111
$ws.Range("K102").Value = 2293.111
$ws.Range("M102").Value = -671.1109999999999

$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("M110").ClearContents()
$ws.Range("N110").ClearContents()

$ws.Range("H122").Value = 1827.2222
$ws.Range("I122").Value = 2162.5
$ws.Range("J122").Value = 1731.4286
$ws.Range("K122").Value = 6487.5
$ws.Range("L122").Value = 5194.2858
$ws.Range("M122").Value = -4037.5
$ws.Range("N122").Value = -10094.2858

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1954.1666
$ws.Range("I99").Value = 1578.5714
$ws.Range("J99").Value = 2108.8235
$ws.Range("K99").Value = 1578.5714
$ws.Range("L99").Value = 2108.8235
$ws.Range("M99").Value = -80.57140000000004
$ws.Range("N99").Value = -5104.8235

$ws.Range("H103").Value = 30539.375
$ws.Range("J103").Value = 30539.375
$ws.Range("L103").Value = 30539.375
$ws.Range("N103").Value = -32883.375

$ws.Range("H105").Value = 2398.3333
$ws.Range("I105").Value = 2382.6924
$ws.Range("K105").Value = 2382.6924
$ws.Range("M105").Value = -635.6923999999999

$ws.Range("H107").Value = 2038.4286
$ws.Range("I107").Value = 2261
$ws.Range("K107").Value = 2261
$ws.Range("M107").Value = -341

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 810
$ws.Range("I16").Value = 583.3333
$ws.Range("J16").Value = 1150
$ws.Range("K16").Value = 583.3333
$ws.Range("L16").Value = 1150
$ws.Range("M16").Value = -296.3333
$ws.Range("N16").Value = -1724

$ws.Range("H31").Value = 1416.7288
$ws.Range("I31").Value = 848.34375
$ws.Range("J31").Value = 2090.3704
$ws.Range("K31").Value = 848.34375
$ws.Range("L31").Value = 2090.3704
$ws.Range("M31").Value = -553.34375
$ws.Range("N31").Value = -2680.3704

$ws.Range("H34").Value = 1416.7288
$ws.Range("I34").Value = 848.34375
$ws.Range("J34").Value = 2090.3704
$ws.Range("K34").Value = 848.34375
$ws.Range("L34").Value = 2090.3704
$ws.Range("M34").Value = -646.34375
$ws.Range("N34").Value = -2494.3704

$ws.Range("H99").Value = 3625310.2
$ws.Range("J99").Value = 3424.75
$ws.Range("L99").Value = 3424.75
$ws.Range("N99").Value = -6420.75

$ws.Range("H105").Value = 587.05884
$ws.Range("J105").Value = 580
$ws.Range("L105").Value = 580
$ws.Range("N105").Value = -4074

$ws.Range("H107").Value = 1865.8125
$ws.Range("I107").Value = 843.75
$ws.Range("J107").Value = 2887.875
$ws.Range("K107").Value = 843.75
$ws.Range("L107").Value = 2887.875
$ws.Range("M107").Value = 1076.25
$ws.Range("N107").Value = -6727.875

$ws.Range("H113").Value = 810
$ws.Range("I113").Value = 583.3333
$ws.Range("J113").Value = 1150
$ws.Range("K113").Value = 583.3333
$ws.Range("L113").Value = 1150
$ws.Range("M113").Value = 1586.6667
$ws.Range("N113").Value = -5490

$ws.Range("H126").Value = 3625310.2
$ws.Range("J126").Value = 3424.75
$ws.Range("L126").Value = 10274.25
$ws.Range("N126").Value = -15214.25

$ws.Range("H132").Value = 599458.2
$ws.Range("I132").Value = 1917.098
$ws.Range("K132").Value = 5751.294
$ws.Range("M132").Value = -3221.294

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 842185.7
$ws.Range("I113").Value = 2020666.5
$ws.Range("J113").Value = 413.66666
$ws.Range("K113").Value = 6061999.5
$ws.Range("L113").Value = 1240.99998
$ws.Range("M113").Value = -6059829.5
$ws.Range("N113").Value = -5580.999980000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 24000
$ws.Range("J103").Value = 24000
$ws.Range("L103").Value = 24000
$ws.Range("N103").Value = -26344

$ws.Range("H107").Value = 1390.7142
$ws.Range("I107").Value = 1347.5
$ws.Range("J107").Value = 1408
$ws.Range("K107").Value = 1347.5
$ws.Range("L107").Value = 1408
$ws.Range("M107").Value = 572.5
$ws.Range("N107").Value = -5248

$ws.Range("H113").Value = 2137
$ws.Range("I113").Value = 1911
$ws.Range("J113").Value = 2250
$ws.Range("K113").Value = 1911
$ws.Range("L113").Value = 2250
$ws.Range("M113").Value = 259
$ws.Range("N113").Value = -6590

$ws.Range("H132").Value = 2780677.8
$ws.Range("I132").Value = 3219.3809
$ws.Range("J132").Value = 6669119.5
$ws.Range("K132").Value = 9658.1427
$ws.Range("L132").Value = 20007358.5
$ws.Range("M132").Value = -7128.1427
$ws.Range("N132").Value = -20012418.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3661.8
$ws.Range("I40").Value = 3501.3333
$ws.Range("J40").Value = 3902.5
$ws.Range("K40").Value = 3501.3333
$ws.Range("L40").Value = 3902.5
$ws.Range("M40").Value = -3365.3333
$ws.Range("N40").Value = -4174.5

$ws.Range("H55").Value = 371.16666
$ws.Range("I55").Value = 126.625
$ws.Range("J55").Value = 460.0909
$ws.Range("K55").Value = 126.625
$ws.Range("L55").Value = 460.0909
$ws.Range("M55").Value = 46.375
$ws.Range("N55").Value = -806.0908999999999

$ws.Range("H61").Value = 3467.0908
$ws.Range("I61").Value = 1898.8
$ws.Range("J61").Value = 4774
$ws.Range("K61").Value = 1898.8
$ws.Range("L61").Value = 4774
$ws.Range("M61").Value = -1696.8
$ws.Range("N61").Value = -5178

$ws.Range("H100").Value = 1750.3334
$ws.Range("I100").Value = 1002
$ws.Range("J100").Value = 1900
$ws.Range("K100").Value = 1002
$ws.Range("L100").Value = 1900
$ws.Range("M100").Value = -461
$ws.Range("N100").Value = -2982

$ws.Range("H113").Value = 3467.0908
$ws.Range("I113").Value = 1898.8
$ws.Range("J113").Value = 4774
$ws.Range("K113").Value = 1898.8
$ws.Range("L113").Value = 4774
$ws.Range("M113").Value = 271.2
$ws.Range("N113").Value = -9114

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 16896
$ws.Range("I107").Value = 16896
$ws.Range("K107").Value = 50688
$ws.Range("M107").Value = -48768

$ws.Range("H113").Value = 617.4545000000001
$ws.Range("I113").Value = 384
$ws.Range("J113").Value = 897.6
$ws.Range("K113").Value = 1152
$ws.Range("L113").Value = 2692.8
$ws.Range("M113").Value = 1018
$ws.Range("N113").Value = -7032.8

$ws.Range("H126").Value = 1700.2858
$ws.Range("I126").Value = 1567.3334
$ws.Range("J126").Value = 1800
$ws.Range("K126").Value = 4702.0002
$ws.Range("L126").Value = 5400
$ws.Range("M126").Value = -2232.0002
$ws.Range("N126").Value = -10340

$ws.Range("H136").Value = 4289.3335
$ws.Range("I136").Value = 4606.2
$ws.Range("K136").Value = 13818.6
$ws.Range("M136").Value = -11268.6
